# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.810.72'
$ws.Range("E2").Value = '  -1.55%  '
$ws.Range("D3").Value = '3.501.32'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.97'
$ws.Range("E5").Value = '  -1.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.84'
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("D7").Value = '3.504.64'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.484'
$ws.Range("E9").Value = '  -1.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.124'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.16'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("D13").Value = '4.079.59'
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.38'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.118'
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.490.13'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000177'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '63.816.71'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.90'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.08'
$ws.Range("E20").Value = '  -1.63%  '
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '383.36'
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.577'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").Value = '3.635.45'
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.37'
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("E27").Value = '  +3.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.59'
$ws.Range("E28").Value = '  +0.55%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.52'
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.28'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("D33").Value = '3.505.74'
$ws.Range("E33").Value = '  -0.40%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.48'
$ws.Range("E35").Value = '  -2.87%  '
$ws.Range("E36").Value = '  +0.88%  '
$ws.Range("E37").Value = '  +2.72%  '
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '160.54'
$ws.Range("E40").Value = '  -4.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0793'
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.68'
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.812'
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.72'
$ws.Range("E45").Value = '  -2.37%  '
$ws.Range("E46").Value = '  -2.23%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  -1.44%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.429.35'
$ws.Range("E49").Value = '  +2.39%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.84'
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.900'
$ws.Range("E51").Value = '  +1.40%  '
